# ============================================================================
# Commit: "Add files via upload"
#
# The document was attached to a SharePoint document library. Saving it
# there stamps three "custom XML part" payloads into the package:
#   customXml/item1.xml - the library's content-type schema (managed-metadata
#                          columns: media service fields, shared-with,
#                          taxonomy catch-all, etc.) plus the InfoPath
#                          PartnerControls person/term helper schema.
#   customXml/item2.xml - the document-library form-template pointer (which
#                          ASPX forms render the Display/Edit/New views).
#   customXml/item3.xml - this document's own managed-metadata values (an
#                          empty taxonomy "Image Tags" term set and a nil
#                          TaxCatchAll column).
# Each item gets its own datastore GUID (itemProps1-3.xml) minted by Word
# when the part is created. None of this touches the visible Word content
# (body text, styles, sections, ...) - it is pure package/custom-XML
# metadata, so the whole edit is expressed purely through
# Document.CustomXMLParts.Add, once per part, in upload order.
# ============================================================================

$d = $word.ActiveDocument

# --- customXml/item1.xml : SharePoint content-type schema ------------------
$item1Xml = @'
<?xml version="1.0" encoding="utf-8"?>
<ct:contentTypeSchema xmlns:ct="http://schemas.microsoft.com/office/2006/metadata/contentType" xmlns:ma="http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes" ct:_="" ma:_="" ma:contentTypeName="Document" ma:contentTypeID="0x010100A0B2487FFFB16146B6CCE0AD83A58ADF" ma:contentTypeVersion="14" ma:contentTypeDescription="Create a new document." ma:contentTypeScope="" ma:versionID="6389e2a6b34849e264c533fb417dfcbe">
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:ns2="c4d1b88e-37a3-489e-92e6-75abb993d69e" xmlns:ns3="e65f5a12-863a-457f-ad99-6387d0246a86" targetNamespace="http://schemas.microsoft.com/office/2006/metadata/properties" ma:root="true" ma:fieldsID="22460c8de2afaaab4523aa0dfde1c706" ns2:_="" ns3:_="">
    <xsd:import namespace="c4d1b88e-37a3-489e-92e6-75abb993d69e"/>
    <xsd:import namespace="e65f5a12-863a-457f-ad99-6387d0246a86"/>
    <xsd:element name="properties">
      <xsd:complexType>
        <xsd:sequence>
          <xsd:element name="documentManagement">
            <xsd:complexType>
              <xsd:all>
                <xsd:element ref="ns2:MediaServiceMetadata" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceFastMetadata" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceSearchProperties" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceObjectDetectorVersions" minOccurs="0"/>
                <xsd:element ref="ns3:SharedWithUsers" minOccurs="0"/>
                <xsd:element ref="ns3:SharedWithDetails" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceDateTaken" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceGenerationTime" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceEventHashCode" minOccurs="0"/>
                <xsd:element ref="ns2:MediaLengthInSeconds" minOccurs="0"/>
                <xsd:element ref="ns2:lcf76f155ced4ddcb4097134ff3c332f" minOccurs="0"/>
                <xsd:element ref="ns3:TaxCatchAll" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceOCR" minOccurs="0"/>
              </xsd:all>
            </xsd:complexType>
          </xsd:element>
        </xsd:sequence>
      </xsd:complexType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="c4d1b88e-37a3-489e-92e6-75abb993d69e" elementFormDefault="qualified">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="MediaServiceMetadata" ma:index="8" nillable="true" ma:displayName="MediaServiceMetadata" ma:hidden="true" ma:internalName="MediaServiceMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceFastMetadata" ma:index="9" nillable="true" ma:displayName="MediaServiceFastMetadata" ma:hidden="true" ma:internalName="MediaServiceFastMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceSearchProperties" ma:index="10" nillable="true" ma:displayName="MediaServiceSearchProperties" ma:hidden="true" ma:internalName="MediaServiceSearchProperties" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceObjectDetectorVersions" ma:index="11" nillable="true" ma:displayName="MediaServiceObjectDetectorVersions" ma:hidden="true" ma:indexed="true" ma:internalName="MediaServiceObjectDetectorVersions" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceDateTaken" ma:index="14" nillable="true" ma:displayName="MediaServiceDateTaken" ma:hidden="true" ma:indexed="true" ma:internalName="MediaServiceDateTaken" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceGenerationTime" ma:index="15" nillable="true" ma:displayName="MediaServiceGenerationTime" ma:hidden="true" ma:internalName="MediaServiceGenerationTime" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceEventHashCode" ma:index="16" nillable="true" ma:displayName="MediaServiceEventHashCode" ma:hidden="true" ma:internalName="MediaServiceEventHashCode" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaLengthInSeconds" ma:index="17" nillable="true" ma:displayName="MediaLengthInSeconds" ma:hidden="true" ma:internalName="MediaLengthInSeconds" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Unknown"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="lcf76f155ced4ddcb4097134ff3c332f" ma:index="19" nillable="true" ma:taxonomy="true" ma:internalName="lcf76f155ced4ddcb4097134ff3c332f" ma:taxonomyFieldName="MediaServiceImageTags" ma:displayName="Image Tags" ma:readOnly="false" ma:fieldId="{5cf76f15-5ced-4ddc-b409-7134ff3c332f}" ma:taxonomyMulti="true" ma:sspId="831fb091-9b95-43a3-9203-971fe8483b85" ma:termSetId="09814cd3-568e-fe90-9814-8d621ff8fb84" ma:anchorId="fba54fb3-c3e1-fe81-a776-ca4b69148c4d" ma:open="true" ma:isKeyword="false">
      <xsd:complexType>
        <xsd:sequence>
          <xsd:element ref="pc:Terms" minOccurs="0" maxOccurs="1"/>
        </xsd:sequence>
      </xsd:complexType>
    </xsd:element>
    <xsd:element name="MediaServiceOCR" ma:index="21" nillable="true" ma:displayName="Extracted Text" ma:internalName="MediaServiceOCR" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note">
          <xsd:maxLength value="255"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="e65f5a12-863a-457f-ad99-6387d0246a86" elementFormDefault="qualified">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="SharedWithUsers" ma:index="12" nillable="true" ma:displayName="Shared With" ma:internalName="SharedWithUsers" ma:readOnly="true">
      <xsd:complexType>
        <xsd:complexContent>
          <xsd:extension base="dms:UserMulti">
            <xsd:sequence>
              <xsd:element name="UserInfo" minOccurs="0" maxOccurs="unbounded">
                <xsd:complexType>
                  <xsd:sequence>
                    <xsd:element name="DisplayName" type="xsd:string" minOccurs="0"/>
                    <xsd:element name="AccountId" type="dms:UserId" minOccurs="0" nillable="true"/>
                    <xsd:element name="AccountType" type="xsd:string" minOccurs="0"/>
                  </xsd:sequence>
                </xsd:complexType>
              </xsd:element>
            </xsd:sequence>
          </xsd:extension>
        </xsd:complexContent>
      </xsd:complexType>
    </xsd:element>
    <xsd:element name="SharedWithDetails" ma:index="13" nillable="true" ma:displayName="Shared With Details" ma:internalName="SharedWithDetails" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note">
          <xsd:maxLength value="255"/>
        </xsd:restriction>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="TaxCatchAll" ma:index="20" nillable="true" ma:displayName="Taxonomy Catch All Column" ma:hidden="true" ma:list="{d17c2052-fd1d-4aab-90cb-88c44f3c2809}" ma:internalName="TaxCatchAll" ma:showField="CatchAllData" ma:web="e65f5a12-863a-457f-ad99-6387d0246a86">
      <xsd:complexType>
        <xsd:complexContent>
          <xsd:extension base="dms:MultiChoiceLookup">
            <xsd:sequence>
              <xsd:element name="Value" type="dms:Lookup" maxOccurs="unbounded" minOccurs="0" nillable="true"/>
            </xsd:sequence>
          </xsd:extension>
        </xsd:complexContent>
      </xsd:complexType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:dc="http://purl.org/dc/elements/1.1/" xmlns:dcterms="http://purl.org/dc/terms/" xmlns:odoc="http://schemas.microsoft.com/internal/obd" targetNamespace="http://schemas.openxmlformats.org/package/2006/metadata/core-properties" elementFormDefault="qualified" attributeFormDefault="unqualified" blockDefault="#all">
    <xsd:import namespace="http://purl.org/dc/elements/1.1/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dc.xsd"/>
    <xsd:import namespace="http://purl.org/dc/terms/" schemaLocation="http://dublincore.org/schemas/xmls/qdc/2003/04/02/dcterms.xsd"/>
    <xsd:element name="coreProperties" type="CT_coreProperties"/>
    <xsd:complexType name="CT_coreProperties">
      <xsd:all>
        <xsd:element ref="dc:creator" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dcterms:created" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dc:identifier" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="contentType" minOccurs="0" maxOccurs="1" type="xsd:string" ma:index="0" ma:displayName="Content Type"/>
        <xsd:element ref="dc:title" minOccurs="0" maxOccurs="1" ma:index="4" ma:displayName="Title"/>
        <xsd:element ref="dc:subject" minOccurs="0" maxOccurs="1"/>
        <xsd:element ref="dc:description" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="keywords" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element ref="dc:language" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="category" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element name="version" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element name="revision" minOccurs="0" maxOccurs="1" type="xsd:string">
          <xsd:annotation>
            <xsd:documentation>
                        This value indicates the number of saves or revisions. The application is responsible for updating this value after each revision.
                    </xsd:documentation>
          </xsd:annotation>
        </xsd:element>
        <xsd:element name="lastModifiedBy" minOccurs="0" maxOccurs="1" type="xsd:string"/>
        <xsd:element ref="dcterms:modified" minOccurs="0" maxOccurs="1"/>
        <xsd:element name="contentStatus" minOccurs="0" maxOccurs="1" type="xsd:string"/>
      </xsd:all>
    </xsd:complexType>
  </xsd:schema>
  <xs:schema xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" xmlns:xs="http://www.w3.org/2001/XMLSchema" targetNamespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" elementFormDefault="qualified" attributeFormDefault="unqualified">
    <xs:element name="Person">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:DisplayName" minOccurs="0"/>
          <xs:element ref="pc:AccountId" minOccurs="0"/>
          <xs:element ref="pc:AccountType" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="DisplayName" type="xs:string"/>
    <xs:element name="AccountId" type="xs:string"/>
    <xs:element name="AccountType" type="xs:string"/>
    <xs:element name="BDCAssociatedEntity">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:BDCEntity" minOccurs="0" maxOccurs="unbounded"/>
        </xs:sequence>
        <xs:attribute ref="pc:EntityNamespace"/>
        <xs:attribute ref="pc:EntityName"/>
        <xs:attribute ref="pc:SystemInstanceName"/>
        <xs:attribute ref="pc:AssociationName"/>
      </xs:complexType>
    </xs:element>
    <xs:attribute name="EntityNamespace" type="xs:string"/>
    <xs:attribute name="EntityName" type="xs:string"/>
    <xs:attribute name="SystemInstanceName" type="xs:string"/>
    <xs:attribute name="AssociationName" type="xs:string"/>
    <xs:element name="BDCEntity">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:EntityDisplayName" minOccurs="0"/>
          <xs:element ref="pc:EntityInstanceReference" minOccurs="0"/>
          <xs:element ref="pc:EntityId1" minOccurs="0"/>
          <xs:element ref="pc:EntityId2" minOccurs="0"/>
          <xs:element ref="pc:EntityId3" minOccurs="0"/>
          <xs:element ref="pc:EntityId4" minOccurs="0"/>
          <xs:element ref="pc:EntityId5" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="EntityDisplayName" type="xs:string"/>
    <xs:element name="EntityInstanceReference" type="xs:string"/>
    <xs:element name="EntityId1" type="xs:string"/>
    <xs:element name="EntityId2" type="xs:string"/>
    <xs:element name="EntityId3" type="xs:string"/>
    <xs:element name="EntityId4" type="xs:string"/>
    <xs:element name="EntityId5" type="xs:string"/>
    <xs:element name="Terms">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:TermInfo" minOccurs="0" maxOccurs="unbounded"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="TermInfo">
      <xs:complexType>
        <xs:sequence>
          <xs:element ref="pc:TermName" minOccurs="0"/>
          <xs:element ref="pc:TermId" minOccurs="0"/>
        </xs:sequence>
      </xs:complexType>
    </xs:element>
    <xs:element name="TermName" type="xs:string"/>
    <xs:element name="TermId" type="xs:string"/>
  </xs:schema>
</ct:contentTypeSchema>

'@
try {
    $part1 = $d.CustomXMLParts.Add($item1Xml)
    Write-Output "Added customXml part 1 (content type schema), Id=$($part1.Id)"
} catch {
    Write-Warning "CustomXMLParts.Add failed for item1: $_"
}

# --- customXml/item2.xml : document-library form templates ------------------
$item2Xml = @'
<?xml version="1.0" encoding="utf-8"?>
<?mso-contentType ?>
<FormTemplates xmlns="http://schemas.microsoft.com/sharepoint/v3/contenttype/forms">
  <Display>DocumentLibraryForm</Display>
  <Edit>DocumentLibraryForm</Edit>
  <New>DocumentLibraryForm</New>
</FormTemplates>

'@
try {
    $part2 = $d.CustomXMLParts.Add($item2Xml)
    Write-Output "Added customXml part 2 (form templates), Id=$($part2.Id)"
} catch {
    Write-Warning "CustomXMLParts.Add failed for item2: $_"
}

# --- customXml/item3.xml : this document's managed-metadata values ---------
$item3Xml = @'
<?xml version="1.0" encoding="utf-8"?>
<p:properties xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls">
  <documentManagement>
    <lcf76f155ced4ddcb4097134ff3c332f xmlns="c4d1b88e-37a3-489e-92e6-75abb993d69e">
      <Terms xmlns="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    </lcf76f155ced4ddcb4097134ff3c332f>
    <TaxCatchAll xmlns="e65f5a12-863a-457f-ad99-6387d0246a86" xsi:nil="true"/>
  </documentManagement>
</p:properties>

'@
try {
    $part3 = $d.CustomXMLParts.Add($item3Xml)
    Write-Output "Added customXml part 3 (document properties), Id=$($part3.Id)"
} catch {
    Write-Warning "CustomXMLParts.Add failed for item3: $_"
}

Write-Output "CustomXMLParts.Count = $($d.CustomXMLParts.Count)"
